# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data refresh values to the 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 295.42856
$ws.Range("I33").Value = 194.54546
$ws.Range("K33").Value = 194.54546
$ws.Range("M33").Value = 34.45454000000001

$ws.Range("H112").Value = 1599.9333
$ws.Range("J112").Value = 1683.0385
$ws.Range("L112").Value = 5049.1155
$ws.Range("N112").Value = -7265.1155

$ws.Range("H129").Value = 852.36664
$ws.Range("J129").Value = 866.4815
$ws.Range("L129").Value = 2599.4445
$ws.Range("N129").Value = -12599.4445

$ws.Range("H137").Value = 2688.311
$ws.Range("I137").Value = 2288.1
$ws.Range("J137").Value = 3008.48
$ws.Range("K137").Value = 6864.299999999999
$ws.Range("L137").Value = 9025.440000000001
$ws.Range("M137").Value = -4314.299999999999
$ws.Range("N137").Value = -14125.44

$ws.Range("H138").Value = 2563.5
$ws.Range("I138").Value = 2619.4
$ws.Range("J138").Value = 2556.1448
$ws.Range("K138").Value = 7858.200000000001
$ws.Range("L138").Value = 7668.4344
$ws.Range("M138").Value = -2718.200000000001
$ws.Range("N138").Value = -17948.4344

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 257
$ws.Range("I26").Value = 257
$ws.Range("K26").Value = 257
$ws.Range("M26").Value = 73

$ws.Range("H74").Value = 1402.8214
$ws.Range("I74").Value = 990.0952
$ws.Range("K74").Value = 990.0952
$ws.Range("M74").Value = -116.0952

$ws.Range("H77").Value = 1402.8214
$ws.Range("I77").Value = 990.0952
$ws.Range("K77").Value = 4950.476
$ws.Range("M77").Value = -582.4759999999997

$ws.Range("H114").Value = 41999
$ws.Range("J114").Value = 41999
$ws.Range("L114").Value = 41999
$ws.Range("N114").Value = -50677

$ws.Range("H122").Value = 2299.2856
$ws.Range("I122").Value = 1570.25
$ws.Range("J122").Value = 3271.3333
$ws.Range("K122").Value = 4710.75
$ws.Range("L122").Value = 9813.999899999999
$ws.Range("M122").Value = -2260.75
$ws.Range("N122").Value = -14713.9999

$ws.Range("H132").Value = 3957.2
$ws.Range("I132").Value = 3764.842
$ws.Range("J132").Value = 4566.3335
$ws.Range("K132").Value = 11294.526
$ws.Range("L132").Value = 13699.0005
$ws.Range("M132").Value = -8764.526
$ws.Range("N132").Value = -18759.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws.Range("H86").Value = 4166.5
$ws.Range("I86").Value = 4177.857
$ws.Range("K86").Value = 4177.857
$ws.Range("M86").Value = -3054.857

$ws.Range("H89").Value = 4166.5
$ws.Range("I89").Value = 4177.857
$ws.Range("K89").Value = 20889.285
$ws.Range("M89").Value = -15273.285

$ws.Range("H94").Value = 8929344
$ws.Range("I94").Value = 11905542
$ws.Range("J94").Value = 751.2857
$ws.Range("K94").Value = 11905542
$ws.Range("L94").Value = 751.2857
$ws.Range("M94").Value = -11905091
$ws.Range("N94").Value = -1653.2857

$ws.Range("H99").Value = 111112420
$ws.Range("I99").Value = 200000960
$ws.Range("J99").Value = 1737.5
$ws.Range("K99").Value = 200000960
$ws.Range("L99").Value = 1737.5
$ws.Range("M99").Value = -199999462
$ws.Range("N99").Value = -4733.5

$ws.Range("H134").Value = 4751.2593
$ws.Range("I134").Value = 1066.1364
$ws.Range("K134").Value = 3198.4092
$ws.Range("M134").Value = -663.4092000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1632.7234
$ws.Range("I31").Value = 1472.2051
$ws.Range("J31").Value = 2415.25
$ws.Range("K31").Value = 1472.2051
$ws.Range("L31").Value = 2415.25
$ws.Range("M31").Value = -1177.2051
$ws.Range("N31").Value = -3005.25

$ws.Range("H34").Value = 1632.7234
$ws.Range("I34").Value = 1472.2051
$ws.Range("J34").Value = 2415.25
$ws.Range("K34").Value = 1472.2051
$ws.Range("L34").Value = 2415.25
$ws.Range("M34").Value = -1270.2051
$ws.Range("N34").Value = -2819.25

$ws.Range("H51").Value = 25000
$ws.Range("J51").Value = 25000
$ws.Range("L51").Value = 25000
$ws.Range("N51").Value = -26472

$ws.Range("H58").Value = 34393
$ws.Range("I58").Value = 1590
$ws.Range("J58").Value = 99999
$ws.Range("K58").Value = 1590
$ws.Range("L58").Value = 99999
$ws.Range("M58").Value = -1387
$ws.Range("N58").Value = -100405

$ws.Range("H59").Value = 27500
$ws.Range("J59").Value = 30000
$ws.Range("L59").Value = 30000
$ws.Range("N59").Value = -32290

$ws.Range("H61").Value = 25000
$ws.Range("J61").Value = 25000
$ws.Range("L61").Value = 25000
$ws.Range("N61").Value = -25696

$ws.Range("H132").Value = 2049.4
$ws.Range("I132").Value = 1706.4736
$ws.Range("J132").Value = 3135.3333
$ws.Range("K132").Value = 5119.4208
$ws.Range("L132").Value = 9405.999899999999
$ws.Range("M132").Value = -2589.4208
$ws.Range("N132").Value = -14465.9999

$ws.Range("H136").Value = 34393
$ws.Range("I136").Value = 1590
$ws.Range("J136").Value = 99999
$ws.Range("K136").Value = 4770
$ws.Range("L136").Value = 299997
$ws.Range("M136").Value = -2220
$ws.Range("N136").Value = -305097

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 8931.5
$ws.Range("J107").Value = 13099.75
$ws.Range("L107").Value = 39299.25
$ws.Range("N107").Value = -43139.25

$ws.Range("H138").Value = 3351.7407
$ws.Range("J138").Value = 2967.2632
$ws.Range("L138").Value = 8901.7896
$ws.Range("N138").Value = -19181.7896

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 2000
$ws.Range("I27").Value = 2000
$ws.Range("K27").Value = 2000
$ws.Range("M27").Value = -1834

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H97").Value = 1054.4615
$ws.Range("I97").Value = 928
$ws.Range("K97").Value = 928
$ws.Range("M97").Value = -432

$ws.Range("H102").Value = 4821.3335
$ws.Range("I102").Value = 4000
$ws.Range("K102").Value = 4000
$ws.Range("M102").Value = -2378

$ws.Range("H132").Value = 5879.3125
$ws.Range("I132").Value = 7878.8335
$ws.Range("K132").Value = 23636.5005
$ws.Range("M132").Value = -21106.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2131.2778
$ws.Range("I7").Value = 1787.875
$ws.Range("J7").Value = 2406
$ws.Range("K7").Value = 1787.875
$ws.Range("L7").Value = 2406
$ws.Range("M7").Value = -1675.875
$ws.Range("N7").Value = -2630

$ws.Range("H126").Value = 2131.2778
$ws.Range("I126").Value = 1787.875
$ws.Range("J126").Value = 2406
$ws.Range("K126").Value = 5363.625
$ws.Range("L126").Value = 7218
$ws.Range("M126").Value = -2893.625
$ws.Range("N126").Value = -12158

$ws.Range("H132").Value = 2480.5134
$ws.Range("I132").Value = 2199.3333
$ws.Range("J132").Value = 2849.5625
$ws.Range("K132").Value = 6597.999899999999
$ws.Range("L132").Value = 8548.6875
$ws.Range("M132").Value = -4067.999899999999
$ws.Range("N132").Value = -13608.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 48310136
$ws.Range("I126").Value = 65360200
$ws.Range("J126").Value = 1620
$ws.Range("K126").Value = 196080600
$ws.Range("L126").Value = 4860
$ws.Range("M126").Value = -196078130
$ws.Range("N126").Value = -9800

$ws.Range("H132").Value = 1869.193
$ws.Range("I132").Value = 1716.82
$ws.Range("J132").Value = 2957.5715
$ws.Range("K132").Value = 5150.46
$ws.Range("L132").Value = 8872.7145
$ws.Range("M132").Value = -2620.46
$ws.Range("N132").Value = -13932.7145
